# Apply "cost parametes for h2 tech in summary" edit.
#
# Summary_old: update HTSE FOM ($/MW-year) / VOM ($/MWh) cells (F2:G2),
# converting them from aggregate totals to per-MW figures.
#
# Summary: add four new columns (CAPEX $/MWe, FOM $/MWe-year, VOM $/MWhe,
# Life (y)) carrying per-technology cost parameters pulled over from
# Summary_old, for every data row (HTSE / PEM / Alkaline).

$wb = $excel.ActiveWorkbook

# ---- Summary_old -----------------------------------------------------
$wsOld = $wb.Worksheets.Item("Summary_old")

$wsOld.Range("F2").Value = 30313.200000000001
$wsOld.Range("G2").Value = 2.76

$wsOld.Range("E2:G2").Select() | Out-Null

# ---- Summary -----------------------------------------------------------
$wsSum = $wb.Worksheets.Item("Summary")

$wsSum.Range("G1").Value = "CAPEX ($/MWe)"
$wsSum.Range("H1").Value = "FOM ($/MWe-year)"
$wsSum.Range("I1").Value = "VOM ($/MWhe)"
$wsSum.Range("J1").Value = "Life (y)"

# HTSE rows (2-6) <- Summary_old row 2 (HTSE)
$wsSum.Range("G2:G6").Value = 646487
$wsSum.Range("H2:H6").Value = 30313.200000000001
$wsSum.Range("I2:I6").Value = 2.76
$wsSum.Range("J2:J6").Value = 20

# PEM rows (7-11) <- Summary_old row 5 (PEM)
$wsSum.Range("G7:G11").Value = 1500000
$wsSum.Range("H7:H11").Value = 12800
$wsSum.Range("I7:I11").Value = 1.3
$wsSum.Range("J7:J11").Value = 20

# Alkaline rows (12-16) <- Summary_old row 4 (Alkaline)
$wsSum.Range("G12:G16").Value = 743865
$wsSum.Range("H12:H16").Value = 60020
$wsSum.Range("I12:I16").Value = 0
$wsSum.Range("J12:J16").Value = 20

$wsSum.Columns.Item(7).ColumnWidth = 12.67
$wsSum.Columns.Item(8).ColumnWidth = 15.83
$wsSum.Columns.Item(9).ColumnWidth = 12.67

$wsSum.Range("L5").Select() | Out-Null

$wsSum.Activate() | Out-Null
